# Working IEEE reader API.
# Applies the shape/content changes captured in the commit diff:
#  - "loads" sheet gains v_nom_kv / s_base_mva (inserted before the old
#    v_nom_pu column) and g_shunt_pu / b_shunt_pu (appended at the end).
#  - "trafos" sheet gains six new trailing headers: idx_hv, idx_lv, tap_pos,
#    tap_change, tap_min, tap_max.
#  - Sheet selections / the active sheet are updated to match the saved
#    view state in the new workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "loads" sheet: widen the table with new columns, preserving the
#    existing column-D width formatting (so it keeps describing the cell
#    that now holds v_nom_pu, matching the target file).
# ---------------------------------------------------------------------
$wsLoads = $wb.Worksheets.Item("loads")

$oldB1 = $wsLoads.Range("B1").Value()
$oldC1 = $wsLoads.Range("C1").Value()
$oldD1 = $wsLoads.Range("D1").Value()
$oldE1 = $wsLoads.Range("E1").Value()

$oldB2 = $wsLoads.Range("B2").Value()
$oldC2 = $wsLoads.Range("C2").Value()
$oldD2 = $wsLoads.Range("D2").Value()
$oldE2 = $wsLoads.Range("E2").Value()

# Move the previous v_nom_pu / p_nom_mw / q_nom_mvar / bus_idx columns two
# slots to the right (D,E,F,G) to make room for the two new leading columns.
$wsLoads.Range("D1").Value = $oldB1
$wsLoads.Range("E1").Value = $oldC1
$wsLoads.Range("F1").Value = $oldD1
$wsLoads.Range("G1").Value = $oldE1

$wsLoads.Range("D2").Value = $oldB2
$wsLoads.Range("E2").Value = $oldC2
$wsLoads.Range("F2").Value = $oldD2
$wsLoads.Range("G2").Value = $oldE2

# New leading columns.
$wsLoads.Range("B1").Value = "v_nom_kv"
$wsLoads.Range("C1").Value = "s_base_mva"
$wsLoads.Range("B2").Value = 22
$wsLoads.Range("C2").Value = 100

# New trailing columns.
$wsLoads.Range("H1").Value = "g_shunt_pu"
$wsLoads.Range("I1").Value = "b_shunt_pu"
$wsLoads.Range("H2").Value = 0
$wsLoads.Range("I2").Value = 0.00484

# ---------------------------------------------------------------------
# 2. "trafos" sheet: append six new headers (no data rows yet).
# ---------------------------------------------------------------------
$wsTrafos = $wb.Worksheets.Item("trafos")

$wsTrafos.Range("I1").Value = "idx_hv"
$wsTrafos.Range("J1").Value = "idx_lv"
$wsTrafos.Range("K1").Value = "tap_pos"
$wsTrafos.Range("L1").Value = "tap_change"
$wsTrafos.Range("M1").Value = "tap_min"
$wsTrafos.Range("N1").Value = "tap_max"

# ---------------------------------------------------------------------
# 3. View state: update each sheet's selection and which tab is active.
# ---------------------------------------------------------------------
$wsLines = $wb.Worksheets.Item("lines")
$wsLines.Range("B1").Select()

$wsGens = $wb.Worksheets.Item("gens")
$wsGens.Range("A1:F1").Select()

$wsTrafos.Range("A1:N1").Select()

# "loads" is the active tab in the edited workbook, and its own selection
# moved to B1 - activate it last so it ends up the saved ActiveSheet.
$wsLoads.Range("B1").Select()
$wsLoads.Activate()
